# Update EUR->ARS rate: 2025-10-05T15:18:05Z
# Appends a new row to the quote history sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 59

# Column A holds a literal date-like string ("2025-10-05"). Force the cell's
# number format to Text first so Excel stores it as the exact string instead
# of auto-converting it into a serial date value.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-10-05"

$ws.Cells.Item($row, 2).Value = "15:18:05"
$ws.Cells.Item($row, 3).Value = "1.00 EUR = 1,794.1737"
